$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Cours 6): replace the old "Mongoose" topic with the new markdown
# links for the "Mongoose" lesson and its exercise.
$ws.Range("C7").Value = "[Mongoose](introduction_mongoose.md)"
$ws.Range("D7").Value = "[Exercice 6 - Mongoose](exercice6_mongoose.md)"

# Row 8 (Cours 7): replace the old "Mongoose" topic with the new markdown
# links for the "Mongoose - la suite" lesson and its exercise.
$ws.Range("C8").Value = "[Mongoose - la suite](mongoose2.md)"
$ws.Range("D8").Value = "[Exercice 7 - Mongoose](exercice7_mongoose.md)"

# Update the active selection to match the saved workbook state.
$ws.Range("D10").Select()
